# Auto-generated script to apply the 2024-10-18 YTD data update
# to cta-violent-crime-ytd.xlsx across the Citywide Totals, By Neighborhood,
# and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 93
$ws.Range("K2").Value = 115
$ws.Range("B3").Value = 65
$ws.Range("E3").Value = 108
$ws.Range("J3").Value = 175
$ws.Range("B4").Value = 10
$ws.Range("B6").Value = 306
$ws.Range("C6").Value = 379
$ws.Range("D6").Value = 329
$ws.Range("E6").Value = 351
$ws.Range("F6").Value = 412
$ws.Range("H6").Value = 369
$ws.Range("I6").Value = 415
$ws.Range("J6").Value = 331
$ws.Range("B7").Value = 414
$ws.Range("C7").Value = 509
$ws.Range("D7").Value = 517
$ws.Range("E7").Value = 525
$ws.Range("F7").Value = 587
$ws.Range("H7").Value = 579
$ws.Range("I7").Value = 697
$ws.Range("J7").Value = 616
$ws.Range("K7").Value = 720

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("C6").Value = 27
$ws.Range("C7").Value = 30

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J3").Value = 9
$ws.Range("C6").Value = 26
$ws.Range("C7").Value = 31
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("J5").Value = 15
$ws.Range("J6").Value = 20
$ws.Range("B4").Value = 1

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B7").Value = 31

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J5").Value = 20
$ws.Range("I6").Value = 4
$ws.Range("J7").Value = 6
$ws.Range("F19").Value = 19
$ws.Range("K19").Value = 28
$ws.Range("B21").Value = 6
$ws.Range("B28").Value = 31
$ws.Range("E29").Value = 6
$ws.Range("C32").Value = 30
$ws.Range("B35").Value = 7
$ws.Range("C36").Value = 31
$ws.Range("J36").Value = 33
$ws.Range("J43").Value = 4
$ws.Range("D53").Value = 63
$ws.Range("J53").Value = 96
$ws.Range("H62").Value = 7
$ws.Range("I77").Value = 40
$ws.Range("E94").Value = 7
$ws.Range("I94").Value = 9
$ws.Range("B97").Value = 4
$ws.Range("B98").Value = 414
$ws.Range("C98").Value = 509
$ws.Range("D98").Value = 517
$ws.Range("E98").Value = 525
$ws.Range("F98").Value = 587
$ws.Range("H98").Value = 579
$ws.Range("I98").Value = 697
$ws.Range("J98").Value = 616
$ws.Range("K98").Value = 720

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("D6").Value = 37
$ws.Range("J6").Value = 51
$ws.Range("D7").Value = 63
$ws.Range("J7").Value = 96

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("B5").Value = 5
$ws.Range("B6").Value = 7

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("H6").Value = 6
$ws.Range("H7").Value = 7

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("E3").Value = 3
$ws.Range("E6").Value = 6

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I6").Value = 26
$ws.Range("I7").Value = 40

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("B3").Value = 1

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("B6").Value = 4

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 6
$ws.Range("F6").Value = 12
$ws.Range("F7").Value = 19
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("J3").Value = 1

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("J6").Value = 6

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J2").Value = 3

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J6").Value = 4

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("E5").Value = 6
$ws.Range("I5").Value = 7
$ws.Range("E6").Value = 7
$ws.Range("I6").Value = 9

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 4

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("B6").Value = 5
$ws.Range("B7").Value = 6
